$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 311
$ws1.Range("F4").Value = 8145
$ws1.Range("F5").Value = 5936
$ws1.Range("F6").Value = 502
$ws1.Range("G9").Value = "不可售"
$ws1.Range("F10").Value = 301
$ws1.Range("F11").Value = 609
$ws1.Range("F12").Value = 71

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 311
$ws4.Range("F4").Value = 8145
$ws4.Range("F5").Value = 5937
$ws4.Range("F6").Value = 502
$ws4.Range("G9").Value = "不可售"
$ws4.Range("F10").Value = 301
$ws4.Range("F15").Value = 609
$ws4.Range("F16").Value = 71
